# Update the "dSF" column (F) values for several rows on Sheet1.
# These edits correspond to a repull / recalculation of the dSF data series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -1
    4  = 0
    5  = -10
    7  = -7
    9  = -7
    10 = -3
    11 = -2
    13 = -6
    14 = -4
    15 = 16
    16 = -1
    17 = -4
    18 = -9
    21 = -9
    24 = -4
    25 = -5
    27 = -4
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
